# Customer Upload.xlsx revision
# - Replace sample "First Name"/"Last Name" values (John Kennedy -> Juan Dela Cruz)
# - Add two new tracking columns: "trackingurn" (S) and "source" (T), with sample
#   values "454dfdfasd34343" and "BRK"
# - Update the sheet's active selection to reflect where the editor was working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sample data values (order chosen so the shared-string table matches the
# target layout: source, 454dfdfasd34343, BRK, Juan, Dela Cruz, trackingurn)
$ws.Range("T1").Value = "source"
$ws.Range("S2").Value = "454dfdfasd34343"
$ws.Range("T2").Value = "BRK"
$ws.Range("C2").Value = "Juan"
$ws.Range("D2").Value = "Dela Cruz"
$ws.Range("S1").Value = "trackingurn"

# Match the bold/filled header formatting used by the other header cells
$ws.Range("R1").Copy() | Out-Null
$ws.Range("S1:T1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Give the new "trackingurn" column a sensible width similar to the other
# text columns (bestfit-style sizing)
$ws.Columns.Item(19).ColumnWidth = 15.04

# Leave the selection where the editor last left it
$ws.Range("P6").Select() | Out-Null
